# Auto-generated PowerShell Excel COM-interop script
# Updates the 'want-to-go count' (column F) on each worksheet to match
# the latest generated snapshot of the source data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1697
$ws.Range("F3").Value = 9753
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = 827
$ws.Range("F6").Value = 643
$ws.Range("F7").Value = 274
$ws.Range("F9").Value = 86
$ws.Range("F11").Value = 1463
$ws.Range("F12").Value = 583
$ws.Range("F13").Value = 65
$ws.Range("F14").Value = 1529
$ws.Range("F15").Value = 135
$ws.Range("F16").Value = 334
$ws.Range("F19").Value = 434
$ws.Range("F20").Value = 1131
$ws.Range("F23").Value = 6
$ws.Range("F24").Value = 56
$ws.Range("F25").Value = 305
$ws.Range("F27").Value = 282
$ws.Range("F30").Value = 654
$ws.Range("F32").Value = 14
$ws.Range("F33").Value = 191
$ws.Range("F34").Value = 92
$ws.Range("F35").Value = 140
$ws.Range("F39").Value = 407
$ws.Range("F40").Value = 659
$ws.Range("F42").Value = 761
$ws.Range("F43").Value = 343
$ws.Range("F44").Value = 299
$ws.Range("F45").Value = 335
$ws.Range("F47").Value = 335

$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 14
$ws.Range("F11").Value = 710
$ws.Range("F21").Value = 1007
$ws.Range("F23").Value = 317
$ws.Range("F25").Value = 302
$ws.Range("F29").Value = 13
$ws.Range("F31").Value = 344
$ws.Range("F34").Value = 190
$ws.Range("F37").Value = 143
$ws.Range("F39").Value = 37
$ws.Range("F44").Value = 48

$ws = $wb.Worksheets.Item(3)
$ws.Range("F6").Value = 2450
$ws.Range("F7").Value = 3885
$ws.Range("F8").Value = 37
$ws.Range("F10").Value = 182
$ws.Range("F11").Value = 153

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1697
$ws.Range("F4").Value = 9754
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 3885
$ws.Range("F7").Value = 182
$ws.Range("F8").Value = 182
$ws.Range("F10").Value = 1463
$ws.Range("F11").Value = 583
$ws.Range("F12").Value = 1529
$ws.Range("F13").Value = 135
$ws.Range("F14").Value = 334
$ws.Range("F17").Value = 434
$ws.Range("F18").Value = 1131
$ws.Range("F22").Value = 6
$ws.Range("F23").Value = 56
$ws.Range("F24").Value = 1007
$ws.Range("F25").Value = 305
$ws.Range("F28").Value = 282
$ws.Range("F30").Value = 654
$ws.Range("F33").Value = 14
$ws.Range("F34").Value = 191
$ws.Range("F35").Value = 344
$ws.Range("F40").Value = 408
$ws.Range("F41").Value = 190
$ws.Range("F42").Value = 659
$ws.Range("F44").Value = 761
$ws.Range("F45").Value = 343
$ws.Range("F46").Value = 37
$ws.Range("F47").Value = 299
$ws.Range("F48").Value = 335
$ws.Range("F49").Value = 335
$ws.Range("F50").Value = 48

